# Generate Report for Handback
# Updates the "latest generated" timestamps recorded on the handback
# status report: the Overview sheet's "Latest HO Xliff Generate Date"
# column, and the per-language sheets' "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-05 01:11:29"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-05 01:11:24"
$zhcn.Range("K2").Value = "2016-09-05 01:11:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-05 01:11:29"
$dede.Range("K2").Value = "2016-09-05 01:11:49"
